# Aggiunto servizio per lanciare in maniera manuale un report
# Update the "Active Report" sheet with the start/end timestamps of the
# last (manual) run, replacing the previous "NO-RUNS" placeholder values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Active Report")

$ws.Range("D2").Value = "10/10/2017 22:00:55"
$ws.Range("E2").Value = "10/10/2017 22:02:21"
